# "Taking Latest to local from Search Module"
# The "Test Cases" sheet has a Runmode column (D) whose values were all
# switched from "N" to "Y" (except the rows that were already "Y"), and the
# sheet view's selection/scroll position moved down to show the bottom of
# the list (D2:D71 selected, scrolled to row 43).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Rows in column D whose Runmode value is currently "N" - flip them to "Y".
$rowsToFlip = @(2,3,4,5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,56,57,59,60,61,62,63,64,65,66,67,68,69,70,71)

foreach ($r in $rowsToFlip) {
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Update the sheet's view so the selection matches the new state: the
# whole Runmode column is selected with D2 as the active cell, and the
# view is scrolled down toward the bottom rows.
$ws.Range("D2:D71").Select()
$excel.ActiveWindow.ScrollRow = 43
